$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell D1 = "stfips" --------------------------------------
# Copy the formatting (bold font + border + centered alignment) already
# used by A1:C1 onto D1 so it reuses the same cell style, then set its
# text value.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D1").Value2 = "stfips"

# --- New data cells D2:D10 = "05" (Census state FIPS code) --------------
# Typing "05" directly would be auto-converted to the number 5, losing the
# leading zero. Build the text in a scratch column via a formula (which
# always yields a text result), then paste-special just the values into
# D2:D10 so each cell ends up holding a plain text/shared-string "05"
# with no special number formatting or quote-prefix styling applied.
$scratch = $ws.Range("F2:F10")
$scratch.Formula = '="05"'
$scratch.Copy() | Out-Null
$ws.Range("D2:D10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null
$excel.CutCopyMode = 0
$scratch.ClearContents()
